$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text so the literal string (with trailing zeros, etc) is preserved.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D17", "D18", "D21", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, matching the authoritative diff.
$ws.Range("D2").Value = "27.948.36"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.867.13"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "312.28"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "0.5030"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("D8").Value = "0.3820"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").Value = "0.08921"
$ws.Range("E9").Value = "  -7.51%  "
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").Value = "41.50"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "6.366"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "20.65"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "1.870.35"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "7.231"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "90.98"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").Value = "6.105"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").Value = "27.974.85"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("D25").Value = "2.272"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.082.45"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.486"
$ws.Range("E27").Value = "  -6.42%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "157.65"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "20.66"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "126.07"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1063"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.052"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.591"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.605"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "9.442"
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06563"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02396"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2180"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.282"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.201"
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "11.49"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "0.6359"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").Value = "4.885"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.18"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5997"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "1.278"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "3.661"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.223"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").Value = "1.988"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "120.48"
$ws.Range("E51").Value = "  -2.47%  "
